# Bond screener "today" reference date moved forward by one day
# (2023-09-29 -> 2023-09-30, serial 45198 -> 45199).
# Column G ("Dni od poprzedniej wyplaty") = today - F ("Data poprzedniej wyplaty")
# Column I ("Dni do nastepnej wyplaty")   = H ("Data nastepnej wyplaty") - today
# Shifting "today" by +1 day increases every G value by 1 and decreases
# every I value by 1, for every row that already has a value in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, 7)   # column G
    if ($gCell.Value2 -ne $null) {
        $gCell.Value2 = $gCell.Value2 + 1
    }

    $iCell = $ws.Cells.Item($r, 9)   # column I
    if ($iCell.Value2 -ne $null) {
        $iCell.Value2 = $iCell.Value2 - 1
    }
}
